# Generate Report for Handback
# Updates the localization-status workbook after a handback transform
# failure: the "Ready for handoff" status becomes "Handback transform
# failed" everywhere it is shown (Overview summary columns + the per
# language "Status" column), and the per-language "Error Detail" cell
# for the affected file is populated with the mismatch diagnostic.
# The "Error Detail" column is also widened to fit the longer message.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row for 4db4fc9d-89cd-49cd-a74a-2930c6c71bac.md -> row 3
# E3 = zh-cn status, F3 = de-de status
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Per-language sheets: row 3 (same file) -> Status column C
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Per-language Error Detail column (P) for that same row
$zhcn.Range("P3").Value = "Handback file name: qbsv1fuk.g10 is different with handoff file name: 4db4fc9d-89cd-49cd-a74a-2930c6c71bac.e29d2024b74e2bcb17f69c0b8b41fe3b230de825.zh-cn."
$dede.Range("P3").Value = "Handback file name: qbsv1fuk.g10 is different with handoff file name: 4db4fc9d-89cd-49cd-a74a-2930c6c71bac.e29d2024b74e2bcb17f69c0b8b41fe3b230de825.de-de."

# Widen the Error Detail column (P / column 16) to fit the new message.
# ColumnWidth is stored with a small constant pixel-padding offset baked
# in by the engine (observed +5/6 character), so back it out here to
# land exactly on a stored width of 40.
$targetWidth = 40 - (5 / 6)
$zhcn.Columns.Item(16).ColumnWidth = $targetWidth
$dede.Columns.Item(16).ColumnWidth = $targetWidth
